$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# Insert four new blank rows at the bottom of the table (147-150). Inserting
# (rather than just writing past the end) makes the new rows inherit the
# number formats (date / text / percentage) already used by column B, C and
# K throughout the table, the same way Excel extends formatting when a user
# types data into the row right below existing formatted data.
$ws.Rows.Item(147).Insert()
$ws.Rows.Item(148).Insert()
$ws.Rows.Item(149).Insert()
$ws.Rows.Item(150).Insert()

# Write the two brand-new date strings first, then the new event name, so
# that the shared-string table allocates them in the same order as the
# source workbook (dates before the event name).
$ws.Range("C149").Value = "2023-11-18"
$ws.Range("C150").Value = "2023-11-19"
$ws.Range("J147").Value = "DOTA 2 ULTRAS DOTA PRO"

# --- Row 147 ---
$ws.Range("A147").Value = 146
$ws.Range("B147").Value = 45247
$ws.Range("C147").Value = "2023-11-17"
$ws.Range("D147").Value = 0
$ws.Range("E147").Value = 1.16
$ws.Range("F147").Formula = "=H146"
$ws.Range("G147").Value = -810
$ws.Range("H147").Formula = "=F147+G147"
$ws.Range("I147").Value = "ESPORTS"
$ws.Range("J147").Value = "DOTA 2 ULTRAS DOTA PRO"
$ws.Range("K147").Formula = "=ROUND((H147/`$F`$31-1)*100, 3)+`$K`$29"

# --- Row 148 ---
$ws.Range("A148").Value = 147
$ws.Range("B148").Value = 45247
$ws.Range("C148").Value = "2023-11-17"
$ws.Range("D148").Value = 0
$ws.Range("E148").Value = 1.16
$ws.Range("F148").Formula = "=H147"
$ws.Range("G148").Value = -200
$ws.Range("H148").Formula = "=F148+G148"
$ws.Range("I148").Value = "ESPORTS"
$ws.Range("J148").Value = "DOTA 2 ULTRAS DOTA PRO"
$ws.Range("K148").Formula = "=ROUND((H148/`$F`$31-1)*100, 3)+`$K`$29"

# --- Row 149 ---
$ws.Range("A149").Value = 148
$ws.Range("B149").Value = 45248
$ws.Range("C149").Value = "2023-11-18"
$ws.Range("D149").Value = 1
$ws.Range("E149").Value = 1.304
$ws.Range("F149").Formula = "=H148"
$ws.Range("G149").Value = 790
$ws.Range("H149").Formula = "=F149+G149"
$ws.Range("I149").Value = "BASKET"
$ws.Range("J149").Value = "NBA"
$ws.Range("K149").Formula = "=ROUND((H149/`$F`$31-1)*100, 3)+`$K`$29"

# --- Row 150 ---
$ws.Range("A150").Value = 149
$ws.Range("B150").Value = 45249
$ws.Range("C150").Value = "2023-11-19"
$ws.Range("D150").Value = 1
$ws.Range("E150").Value = 1.14
$ws.Range("F150").Formula = "=H149"
$ws.Range("G150").Value = 13860
$ws.Range("H150").Formula = "=F150+G150"
$ws.Range("I150").Value = "ESPORTS"
$ws.Range("J150").Value = "WORLDS 2023"
$ws.Range("K150").Formula = "=ROUND((H150/`$F`$31-1)*100, 3)+`$K`$29"

# Reflect the new selection / scroll position left behind by the edit.
$ws.Range("M148").Select()
